# Update countries & provincias Spain
# - Refresh COVID-19 figures for a handful of countries (Suiza, Rumania,
#   Bosnia y Herzegovina, Malta, Brunei, Etiopia, Tanzania, Malaui).
# - The table (A4:H215) is kept sorted by "Casos totales" (col B) descending,
#   so the refreshed totals re-shuffle several rows; update the whole table in
#   place so row positions/labels stay consistent with the new sort order.
# - Bump the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 12:52"

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 367650
$ws.Range("C4").Value = 646
$ws.Range("D4").Value = 19810
$ws.Range("E4").Value = 336897
$ws.Range("F4").Value = 8983
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 10943
$ws.Range("A5").Value = "España"
$ws.Range("B5").Value = 140510
$ws.Range("C5").Value = 3835
$ws.Range("D5").Value = 43208
$ws.Range("E5").Value = 83504
$ws.Range("F5").Value = 7069
$ws.Range("G5").Value = 457
$ws.Range("H5").Value = 13798
$ws.Range("A6").Value = "Italia"
$ws.Range("B6").Value = 132547
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 22837
$ws.Range("E6").Value = 93187
$ws.Range("F6").Value = 3898
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 16523
$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 103375
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 36081
$ws.Range("E7").Value = 65484
$ws.Range("F7").Value = 4895
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1810
$ws.Range("A8").Value = "Francia"
$ws.Range("B8").Value = 98010
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 17250
$ws.Range("E8").Value = 71849
$ws.Range("F8").Value = 7072
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 8911
$ws.Range("A9").Value = "China"
$ws.Range("B9").Value = 81740
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 77167
$ws.Range("E9").Value = 1242
$ws.Range("F9").Value = 211
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 3331
$ws.Range("A10").Value = "Iran"
$ws.Range("B10").Value = 62589
$ws.Range("C10").Value = 2089
$ws.Range("D10").Value = 27039
$ws.Range("E10").Value = 31678
$ws.Range("F10").Value = 3987
$ws.Range("G10").Value = 133
$ws.Range("H10").Value = 3872
$ws.Range("A11").Value = "Reino Unido"
$ws.Range("B11").Value = 51608
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 135
$ws.Range("E11").Value = 46100
$ws.Range("F11").Value = 1559
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 5373
$ws.Range("A12").Value = "Turquia"
$ws.Range("B12").Value = 30217
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 1326
$ws.Range("E12").Value = 28242
$ws.Range("F12").Value = 1415
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 649
$ws.Range("A13").Value = "Belgica"
$ws.Range("B13").Value = 22194
$ws.Range("C13").Value = 1380
$ws.Range("D13").Value = 4157
$ws.Range("E13").Value = 16002
$ws.Range("F13").Value = 1260
$ws.Range("G13").Value = 403
$ws.Range("H13").Value = 2035
$ws.Range("A14").Value = "Suiza"
$ws.Range("B14").Value = 21907
$ws.Range("C14").Value = 250
$ws.Range("D14").Value = 8056
$ws.Range("E14").Value = 13064
$ws.Range("F14").Value = 391
$ws.Range("G14").Value = 22
$ws.Range("H14").Value = 787
$ws.Range("A15").Value = "Paises Bajos"
$ws.Range("B15").Value = 18803
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 250
$ws.Range("E15").Value = 16686
$ws.Range("F15").Value = 1409
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1867
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 16667
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 3616
$ws.Range("E16").Value = 12728
$ws.Range("F16").Value = 426
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 323
$ws.Range("A17").Value = "Austria"
$ws.Range("B17").Value = 12427
$ws.Range("C17").Value = 130
$ws.Range("D17").Value = 4046
$ws.Range("E17").Value = 8138
$ws.Range("F17").Value = 243
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 243
$ws.Range("A18").Value = "Brasil"
$ws.Range("B18").Value = 12240
$ws.Range("C18").Value = 57
$ws.Range("D18").Value = 127
$ws.Range("E18").Value = 11547
$ws.Range("F18").Value = 296
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 566
$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 11730
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 140
$ws.Range("E19").Value = 11279
$ws.Range("F19").Value = 270
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 311
$ws.Range("A20").Value = "Corea del Sur"
$ws.Range("B20").Value = 10331
$ws.Range("C20").Value = 47
$ws.Range("D20").Value = 6694
$ws.Range("E20").Value = 3445
$ws.Range("F20").Value = 55
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 192
$ws.Range("A21").Value = "Israel"
$ws.Range("B21").Value = 9006
$ws.Range("C21").Value = 102
$ws.Range("D21").Value = 683
$ws.Range("E21").Value = 8264
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 59
$ws.Range("A22").Value = "Rusia"
$ws.Range("B22").Value = 7497
$ws.Range("C22").Value = 1154
$ws.Range("D22").Value = 494
$ws.Range("E22").Value = 6945
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 58
$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 7206
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 205
$ws.Range("E23").Value = 6524
$ws.Range("F23").Value = 590
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 477
$ws.Range("A24").Value = "Australia"
$ws.Range("B24").Value = 5908
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 2547
$ws.Range("E24").Value = 3315
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 46
$ws.Range("A25").Value = "Noruega"
$ws.Range("B25").Value = 5866
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 5756
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 78
$ws.Range("A26").Value = "Irlanda"
$ws.Range("B26").Value = 5364
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 5165
$ws.Range("F26").Value = 165
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 174
$ws.Range("A27").Value = "Dinamarca"
$ws.Range("B27").Value = 4978
$ws.Range("C27").Value = 297
$ws.Range("D27").Value = 1378
$ws.Range("E27").Value = 3413
$ws.Range("F27").Value = 144
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 187
$ws.Range("A28").Value = "India"
$ws.Range("B28").Value = 4858
$ws.Range("C28").Value = 80
$ws.Range("D28").Value = 382
$ws.Range("E28").Value = 4339
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 137
$ws.Range("A29").Value = "Chequia"
$ws.Range("B29").Value = 4828
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 127
$ws.Range("E29").Value = 4621
$ws.Range("F29").Value = 86
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 80
$ws.Range("A30").Value = "Chile"
$ws.Range("B30").Value = 4815
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 728
$ws.Range("E30").Value = 4050
$ws.Range("F30").Value = 327
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 37
$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 4532
$ws.Range("C31").Value = 119
$ws.Range("D31").Value = 191
$ws.Range("E31").Value = 4230
$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 111
$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 4417
$ws.Range("C32").Value = 360
$ws.Range("D32").Value = 460
$ws.Range("E32").Value = 3775
$ws.Range("F32").Value = 274
$ws.Range("G32").Value = 6
$ws.Range("H32").Value = 182
$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 4004
$ws.Range("C33").Value = 238
$ws.Range("D33").Value = 429
$ws.Range("E33").Value = 3521
$ws.Range("F33").Value = 28
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 54
$ws.Range("A34").Value = "Malasia"
$ws.Range("B34").Value = 3963
$ws.Range("C34").Value = 170
$ws.Range("D34").Value = 1321
$ws.Range("E34").Value = 2579
$ws.Range("F34").Value = 92
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 63
$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 3906
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 592
$ws.Range("E35").Value = 3222
$ws.Range("F35").Value = 79
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 92
$ws.Range("A36").Value = "Filipinas"
$ws.Range("B36").Value = 3764
$ws.Range("C36").Value = 104
$ws.Range("D36").Value = 84
$ws.Range("E36").Value = 3503
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 177
$ws.Range("A37").Value = "Ecuador"
$ws.Range("B37").Value = 3747
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 100
$ws.Range("E37").Value = 3456
$ws.Range("F37").Value = 156
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 191
$ws.Range("A38").Value = "Luxemburgo"
$ws.Range("B38").Value = 2843
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 500
$ws.Range("E38").Value = 2302
$ws.Range("F38").Value = 33
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 41
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("B39").Value = 2752
$ws.Range("C39").Value = 147
$ws.Range("D39").Value = 551
$ws.Range("E39").Value = 2163
$ws.Range("F39").Value = 41
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 38
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 2738
$ws.Range("C40").Value = 247
$ws.Range("D40").Value = 204
$ws.Range("E40").Value = 2313
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 12
$ws.Range("H40").Value = 221
$ws.Range("A41").Value = "Peru"
$ws.Range("B41").Value = 2561
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 997
$ws.Range("E41").Value = 1472
$ws.Range("F41").Value = 89
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 92
$ws.Range("A42").Value = "Mexico"
$ws.Range("B42").Value = 2439
$ws.Range("C42").Value = 296
$ws.Range("D42").Value = 633
$ws.Range("E42").Value = 1681
$ws.Range("F42").Value = 89
$ws.Range("G42").Value = 31
$ws.Range("H42").Value = 125
$ws.Range("A43").Value = "Finlandia"
$ws.Range("B43").Value = 2308
$ws.Range("C43").Value = 132
$ws.Range("D43").Value = 300
$ws.Range("E43").Value = 1981
$ws.Range("F43").Value = 81
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 27
$ws.Range("A44").Value = "Tailandia"
$ws.Range("B44").Value = 2258
$ws.Range("C44").Value = 38
$ws.Range("D44").Value = 824
$ws.Range("E44").Value = 1407
$ws.Range("F44").Value = 30
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 27
$ws.Range("A45").Value = "Serbia"
$ws.Range("B45").Value = 2200
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 118
$ws.Range("E45").Value = 2024
$ws.Range("F45").Value = 101
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 58
$ws.Range("A46").Value = "Panama"
$ws.Range("B46").Value = 2100
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 14
$ws.Range("E46").Value = 2031
$ws.Range("F46").Value = 88
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 55
$ws.Range("A47").Value = "Emiratos Arabes Unidos"
$ws.Range("B47").Value = 2076
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 167
$ws.Range("E47").Value = 1898
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 11
$ws.Range("A48").Value = "Catar"
$ws.Range("B48").Value = 1832
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 131
$ws.Range("E48").Value = 1697
$ws.Range("F48").Value = 37
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 4
$ws.Range("A49").Value = "Republica Dominicana"
$ws.Range("B49").Value = 1828
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 33
$ws.Range("E49").Value = 1709
$ws.Range("F49").Value = 147
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 86
$ws.Range("A50").Value = "Grecia"
$ws.Range("B50").Value = 1755
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 269
$ws.Range("E50").Value = 1407
$ws.Range("F50").Value = 90
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 79
$ws.Range("A51").Value = "Sudafrica"
$ws.Range("B51").Value = 1686
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 95
$ws.Range("E51").Value = 1579
$ws.Range("F51").Value = 7
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 12
$ws.Range("A52").Value = "Argentina"
$ws.Range("B52").Value = 1628
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 325
$ws.Range("E52").Value = 1250
$ws.Range("F52").Value = 94
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 53
$ws.Range("A53").Value = "Colombia"
$ws.Range("B53").Value = 1579
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 88
$ws.Range("E53").Value = 1445
$ws.Range("F53").Value = 50
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 46
$ws.Range("A54").Value = "Islandia"
$ws.Range("B54").Value = 1562
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 460
$ws.Range("E54").Value = 1096
$ws.Range("F54").Value = 11
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 6
$ws.Range("A55").Value = "Ucrania"
$ws.Range("B55").Value = 1462
$ws.Range("C55").Value = 143
$ws.Range("D55").Value = 28
$ws.Range("E55").Value = 1389
$ws.Range("F55").Value = 16
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 45
$ws.Range("A56").Value = "Argelia"
$ws.Range("B56").Value = 1423
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 90
$ws.Range("E56").Value = 1160
$ws.Range("F56").Value = 46
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 173
$ws.Range("A57").Value = "Singapur"
$ws.Range("B57").Value = 1375
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 344
$ws.Range("E57").Value = 1025
$ws.Range("F57").Value = 25
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 6
$ws.Range("A58").Value = "Egipto"
$ws.Range("B58").Value = 1322
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 259
$ws.Range("E58").Value = 978
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 85
$ws.Range("A59").Value = "Croacia"
$ws.Range("B59").Value = 1222
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 130
$ws.Range("E59").Value = 1076
$ws.Range("F59").Value = 36
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 16
$ws.Range("A60").Value = "Nueva Zelanda"
$ws.Range("B60").Value = 1160
$ws.Range("C60").Value = 54
$ws.Range("D60").Value = 241
$ws.Range("E60").Value = 918
$ws.Range("F60").Value = 14
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1
$ws.Range("A61").Value = "Estonia"
$ws.Range("B61").Value = 1149
$ws.Range("C61").Value = 41
$ws.Range("D61").Value = 69
$ws.Range("E61").Value = 1059
$ws.Range("F61").Value = 12
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 21
$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 1141
$ws.Range("C62").Value = 21
$ws.Range("D62").Value = 88
$ws.Range("E62").Value = 970
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 83
$ws.Range("A63").Value = "Irak"
$ws.Range("B63").Value = 1031
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 344
$ws.Range("E63").Value = 623
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 64
$ws.Range("A64").Value = "Eslovenia"
$ws.Range("B64").Value = 1021
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 102
$ws.Range("E64").Value = 889
$ws.Range("F64").Value = 30
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 30
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 965
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 40
$ws.Range("E65").Value = 904
$ws.Range("F65").Value = 80
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 21
$ws.Range("A66").Value = "Hong Kong"
$ws.Range("B66").Value = 936
$ws.Range("C66").Value = 21
$ws.Range("D66").Value = 236
$ws.Range("E66").Value = 696
$ws.Range("F66").Value = 12
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 4
$ws.Range("A67").Value = "Lituania"
$ws.Range("B67").Value = 880
$ws.Range("C67").Value = 37
$ws.Range("D67").Value = 8
$ws.Range("E67").Value = 857
$ws.Range("F67").Value = 11
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 15
$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 853
$ws.Range("C68").Value = 20
$ws.Range("D68").Value = 87
$ws.Range("E68").Value = 758
$ws.Range("F68").Value = 30
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 8
$ws.Range("A69").Value = "Hungria"
$ws.Range("B69").Value = 817
$ws.Range("C69").Value = 73
$ws.Range("D69").Value = 71
$ws.Range("E69").Value = 699
$ws.Range("F69").Value = 17
$ws.Range("G69").Value = 9
$ws.Range("H69").Value = 47
$ws.Range("A70").Value = "Barein"
$ws.Range("B70").Value = 756
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 458
$ws.Range("E70").Value = 294
$ws.Range("F70").Value = 4
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 4
$ws.Range("A71").Value = "Bosnia y Herzegovina"
$ws.Range("B71").Value = 744
$ws.Range("C71").Value = 70
$ws.Range("D71").Value = 68
$ws.Range("E71").Value = 643
$ws.Range("F71").Value = 4
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 33
$ws.Range("A72").Value = "Kuwait"
$ws.Range("B72").Value = 743
$ws.Range("C72").Value = 78
$ws.Range("D72").Value = 105
$ws.Range("E72").Value = 637
$ws.Range("F72").Value = 23
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 1
$ws.Range("A73").Value = "Crucero"
$ws.Range("B73").Value = 712
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 619
$ws.Range("E73").Value = 82
$ws.Range("F73").Value = 10
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 11
$ws.Range("A74").Value = "Bielorrusia"
$ws.Range("B74").Value = 700
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 53
$ws.Range("E74").Value = 634
$ws.Range("F74").Value = 11
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 13
$ws.Range("A75").Value = "Kazajistan"
$ws.Range("B75").Value = 685
$ws.Range("C75").Value = 23
$ws.Range("D75").Value = 50
$ws.Range("E75").Value = 629
$ws.Range("F75").Value = 16
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 6
$ws.Range("A76").Value = "Camerun"
$ws.Range("B76").Value = 658
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 17
$ws.Range("E76").Value = 632
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 9
$ws.Range("A77").Value = "Azerbaiyan"
$ws.Range("B77").Value = 641
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 44
$ws.Range("E77").Value = 590
$ws.Range("F77").Value = 11
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 7
$ws.Range("A78").Value = "Tunez"
$ws.Range("B78").Value = 596
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 5
$ws.Range("E78").Value = 569
$ws.Range("F78").Value = 39
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 22
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 570
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 30
$ws.Range("E79").Value = 517
$ws.Range("F79").Value = 15
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 23
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("B80").Value = 565
$ws.Range("C80").Value = 16
$ws.Range("D80").Value = 42
$ws.Range("E80").Value = 501
$ws.Range("F80").Value = 26
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 22
$ws.Range("A81").Value = "Letonia"
$ws.Range("B81").Value = 548
$ws.Range("C81").Value = 6
$ws.Range("D81").Value = 16
$ws.Range("E81").Value = 530
$ws.Range("F81").Value = 5
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 2
$ws.Range("A82").Value = "Libano"
$ws.Range("B82").Value = 548
$ws.Range("C82").Value = 7
$ws.Range("D82").Value = 60
$ws.Range("E82").Value = 469
$ws.Range("F82").Value = 27
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 19
$ws.Range("A83").Value = "Eslovaquia"
$ws.Range("B83").Value = 534
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 8
$ws.Range("E83").Value = 524
$ws.Range("F83").Value = 3
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 2
$ws.Range("A84").Value = "Principado de Andorra"
$ws.Range("B84").Value = 525
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 31
$ws.Range("E84").Value = 473
$ws.Range("F84").Value = 12
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 21
$ws.Range("A85").Value = "Uzbekistan"
$ws.Range("B85").Value = 472
$ws.Range("C85").Value = 15
$ws.Range("D85").Value = 30
$ws.Range("E85").Value = 440
$ws.Range("F85").Value = 8
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 2
$ws.Range("A86").Value = "Costa Rica"
$ws.Range("B86").Value = 467
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 18
$ws.Range("E86").Value = 447
$ws.Range("F86").Value = 14
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 2
$ws.Range("A87").Value = "Republica de Chipre"
$ws.Range("B87").Value = 465
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 45
$ws.Range("E87").Value = 411
$ws.Range("F87").Value = 11
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 9
$ws.Range("A88").Value = "Afganistan"
$ws.Range("B88").Value = 423
$ws.Range("C88").Value = 56
$ws.Range("D88").Value = 18
$ws.Range("E88").Value = 394
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 11
$ws.Range("A89").Value = "Uruguay"
$ws.Range("B89").Value = 415
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 123
$ws.Range("E89").Value = 286
$ws.Range("F89").Value = 14
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 6
$ws.Range("A90").Value = "Albania"
$ws.Range("B90").Value = 383
$ws.Range("C90").Value = 6
$ws.Range("D90").Value = 131
$ws.Range("E90").Value = 230
$ws.Range("F90").Value = 7
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 22
$ws.Range("A91").Value = "Taiwan"
$ws.Range("B91").Value = 376
$ws.Range("C91").Value = 3
$ws.Range("D91").Value = 61
$ws.Range("E91").Value = 310
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 5
$ws.Range("A92").Value = "Oman"
$ws.Range("B92").Value = 371
$ws.Range("C92").Value = 40
$ws.Range("D92").Value = 67
$ws.Range("E92").Value = 302
$ws.Range("F92").Value = 3
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 2
$ws.Range("A93").Value = "Burkina Faso"
$ws.Range("B93").Value = 364
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 108
$ws.Range("E93").Value = 238
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 18
$ws.Range("A94").Value = "Cuba"
$ws.Range("B94").Value = 363
$ws.Range("C94").Value = 13
$ws.Range("D94").Value = 18
$ws.Range("E94").Value = 336
$ws.Range("F94").Value = 12
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 9
$ws.Range("A95").Value = "Reunion"
$ws.Range("B95").Value = 349
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 40
$ws.Range("E95").Value = 309
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("A96").Value = "Jordania"
$ws.Range("B96").Value = 349
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 126
$ws.Range("E96").Value = 217
$ws.Range("F96").Value = 5
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 6
$ws.Range("A97").Value = "Costa de Marfil"
$ws.Range("B97").Value = 323
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 41
$ws.Range("E97").Value = 279
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 3
$ws.Range("A98").Value = "Honduras"
$ws.Range("B98").Value = 305
$ws.Range("C98").Value = 7
$ws.Range("D98").Value = 6
$ws.Range("E98").Value = 277
$ws.Range("F98").Value = 10
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 22
$ws.Range("A99").Value = "Malta"
$ws.Range("B99").Value = 293
$ws.Range("C99").Value = 52
$ws.Range("D99").Value = 5
$ws.Range("E99").Value = 288
$ws.Range("F99").Value = 4
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("A100").Value = "Ghana"
$ws.Range("B100").Value = 287
$ws.Range("C100").Value = 73
$ws.Range("D100").Value = 31
$ws.Range("E100").Value = 251
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 5
$ws.Range("A101").Value = "San Marino"
$ws.Range("B101").Value = 277
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 35
$ws.Range("E101").Value = 210
$ws.Range("F101").Value = 14
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 32
$ws.Range("A102").Value = "Estado de Palestina"
$ws.Range("B102").Value = 260
$ws.Range("C102").Value = 6
$ws.Range("D102").Value = 24
$ws.Range("E102").Value = 235
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 1
$ws.Range("A103").Value = "Niger"
$ws.Range("B103").Value = 253
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 26
$ws.Range("E103").Value = 217
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 10
$ws.Range("A104").Value = "Vietnam"
$ws.Range("B104").Value = 245
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 123
$ws.Range("E104").Value = 122
$ws.Range("F104").Value = 8
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("A105").Value = "Mauricio"
$ws.Range("B105").Value = 244
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 7
$ws.Range("E105").Value = 230
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7
$ws.Range("A106").Value = "Montenegro"
$ws.Range("B106").Value = 239
$ws.Range("C106").Value = 6
$ws.Range("D106").Value = 1
$ws.Range("E106").Value = 236
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2
$ws.Range("A107").Value = "Nigeria"
$ws.Range("B107").Value = 238
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 35
$ws.Range("E107").Value = 198
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 5
$ws.Range("A108").Value = "Kirguistan"
$ws.Range("B108").Value = 228
$ws.Range("C108").Value = 12
$ws.Range("D108").Value = 33
$ws.Range("E108").Value = 191
$ws.Range("F108").Value = 5
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 4
$ws.Range("A109").Value = "Senegal"
$ws.Range("B109").Value = 226
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 92
$ws.Range("E109").Value = 132
$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 2
$ws.Range("A110").Value = "Georgia"
$ws.Range("B110").Value = 195
$ws.Range("C110").Value = 7
$ws.Range("D110").Value = 39
$ws.Range("E110").Value = 154
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2
$ws.Range("A111").Value = "Bolivia"
$ws.Range("B111").Value = 194
$ws.Range("C111").Value = 11
$ws.Range("D111").Value = 2
$ws.Range("E111").Value = 178
$ws.Range("F111").Value = 3
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = 14
$ws.Range("A112").Value = "Islas Feroe"
$ws.Range("B112").Value = 184
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 129
$ws.Range("E112").Value = 55
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("A113").Value = "Sri Lanka"
$ws.Range("B113").Value = 180
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 42
$ws.Range("E113").Value = 132
$ws.Range("F113").Value = 5
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 6
$ws.Range("A114").Value = "Venezuela"
$ws.Range("B114").Value = 165
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 65
$ws.Range("E114").Value = 93
$ws.Range("F114").Value = 6
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 7
$ws.Range("A115").Value = "Mayotte"
$ws.Range("B115").Value = 164
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 15
$ws.Range("E115").Value = 147
$ws.Range("F115").Value = 3
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 2
$ws.Range("A116").Value = "Banglades"
$ws.Range("B116").Value = 164
$ws.Range("C116").Value = 41
$ws.Range("D116").Value = 33
$ws.Range("E116").Value = 114
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 5
$ws.Range("H116").Value = 17
$ws.Range("A117").Value = "Consejo Danes para los Refugiados"
$ws.Range("B117").Value = 161
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 5
$ws.Range("E117").Value = 138
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 18
$ws.Range("A118").Value = "Kenia"
$ws.Range("B118").Value = 158
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 4
$ws.Range("E118").Value = 148
$ws.Range("F118").Value = 2
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 6
$ws.Range("A119").Value = "Martinica"
$ws.Range("B119").Value = 151
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 50
$ws.Range("E119").Value = 97
$ws.Range("F119").Value = 20
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 4
$ws.Range("A120").Value = "Isla de Man"
$ws.Range("B120").Value = 150
$ws.Range("C120").Value = 11
$ws.Range("D120").Value = 73
$ws.Range("E120").Value = 76
$ws.Range("F120").Value = 6
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1
$ws.Range("A121").Value = "Guadalupe"
$ws.Range("B121").Value = 139
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 31
$ws.Range("E121").Value = 101
$ws.Range("F121").Value = 14
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 7
$ws.Range("A122").Value = "Brunei"
$ws.Range("B122").Value = 135
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 85
$ws.Range("E122").Value = 49
$ws.Range("F122").Value = 3
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 1
$ws.Range("A123").Value = "Guinea"
$ws.Range("B123").Value = 128
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 5
$ws.Range("E123").Value = 123
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("A124").Value = "Paraguay"
$ws.Range("B124").Value = 115
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 15
$ws.Range("E124").Value = 95
$ws.Range("F124").Value = 1
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 5
$ws.Range("A125").Value = "Camboya"
$ws.Range("B125").Value = 115
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 58
$ws.Range("E125").Value = 57
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("A126").Value = "Gibraltar"
$ws.Range("B126").Value = 109
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 52
$ws.Range("E126").Value = 57
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("A127").Value = "Ruanda"
$ws.Range("B127").Value = 105
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 4
$ws.Range("E127").Value = 101
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("A128").Value = "Trinidad yTobago"
$ws.Range("B128").Value = 105
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 96
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 8
$ws.Range("A129").Value = "Republica de Yibuti"
$ws.Range("B129").Value = 90
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 9
$ws.Range("E129").Value = 81
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("A130").Value = "Madagascar"
$ws.Range("B130").Value = 82
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 2
$ws.Range("E130").Value = 80
$ws.Range("F130").Value = 6
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("A131").Value = "El Salvador"
$ws.Range("B131").Value = 78
$ws.Range("C131").Value = 9
$ws.Range("D131").Value = 5
$ws.Range("E131").Value = 69
$ws.Range("F131").Value = 4
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 4
$ws.Range("A132").Value = "Monaco"
$ws.Range("B132").Value = 77
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 4
$ws.Range("E132").Value = 72
$ws.Range("F132").Value = 4
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1
$ws.Range("A133").Value = "Liechtenstein"
$ws.Range("B133").Value = 77
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 55
$ws.Range("E133").Value = 21
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 1
$ws.Range("A134").Value = "Guatemala"
$ws.Range("B134").Value = 74
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 17
$ws.Range("E134").Value = 54
$ws.Range("F134").Value = 3
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 3
$ws.Range("A135").Value = "Guayana Francesa"
$ws.Range("B135").Value = 72
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 34
$ws.Range("E135").Value = 38
$ws.Range("F135").Value = 1
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("A136").Value = "Aruba"
$ws.Range("B136").Value = 71
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 2
$ws.Range("E136").Value = 69
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("A137").Value = "Barbados"
$ws.Range("B137").Value = 60
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 6
$ws.Range("E137").Value = 52
$ws.Range("F137").Value = 4
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 2
$ws.Range("A138").Value = "Jamaica"
$ws.Range("B138").Value = 59
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 8
$ws.Range("E138").Value = 48
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 3
$ws.Range("A139").Value = "Togo"
$ws.Range("B139").Value = 58
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 23
$ws.Range("E139").Value = 32
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 3
$ws.Range("A140").Value = "Uganda"
$ws.Range("B140").Value = 52
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 52
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("A141").Value = "Etiopia"
$ws.Range("B141").Value = 52
$ws.Range("C141").Value = 8
$ws.Range("D141").Value = 4
$ws.Range("E141").Value = 46
$ws.Range("F141").Value = 1
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 2
$ws.Range("A142").Value = "Mali"
$ws.Range("B142").Value = 47
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 9
$ws.Range("E142").Value = 33
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 5
$ws.Range("A143").Value = "Congo"
$ws.Range("B143").Value = 45
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 2
$ws.Range("E143").Value = 38
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 5
$ws.Range("A144").Value = "Macao"
$ws.Range("B144").Value = 44
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 10
$ws.Range("E144").Value = 34
$ws.Range("F144").Value = 1
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0
$ws.Range("A145").Value = "Polinesia Francesa"
$ws.Range("B145").Value = 42
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 42
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0
$ws.Range("A146").Value = "Islas Caimanes"
$ws.Range("B146").Value = 39
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 1
$ws.Range("E146").Value = 37
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 1
$ws.Range("A147").Value = "Puerto Rico"
$ws.Range("B147").Value = 39
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 1
$ws.Range("E147").Value = 36
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 2
$ws.Range("A148").Value = "Zambia"
$ws.Range("B148").Value = 39
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 5
$ws.Range("E148").Value = 33
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1
$ws.Range("A149").Value = "Bermudas"
$ws.Range("B149").Value = 39
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 17
$ws.Range("E149").Value = 20
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 2
$ws.Range("A150").Value = "San Martin (Parte Holandesa)"
$ws.Range("B150").Value = 37
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 30
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6
$ws.Range("A151").Value = "Guinea-Bisau"
$ws.Range("B151").Value = 33
$ws.Range("C151").Value = 15
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 33
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0
$ws.Range("A152").Value = "Bahamas"
$ws.Range("B152").Value = 33
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 5
$ws.Range("E152").Value = 23
$ws.Range("F152").Value = 1
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 5
$ws.Range("A153").Value = "Guam"
$ws.Range("B153").Value = 32
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 31
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 1
$ws.Range("A154").Value = "San Martin (Parte Francesa)"
$ws.Range("B154").Value = 32
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 7
$ws.Range("E154").Value = 23
$ws.Range("F154").Value = 6
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 2
$ws.Range("A155").Value = "Eritrea"
$ws.Range("B155").Value = 31
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 31
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 0
$ws.Range("A156").Value = "Guyana"
$ws.Range("B156").Value = 31
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 8
$ws.Range("E156").Value = 18
$ws.Range("F156").Value = 8
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 5
$ws.Range("A157").Value = "Gabon"
$ws.Range("B157").Value = 30
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 1
$ws.Range("E157").Value = 28
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 1
$ws.Range("A158").Value = "Benin"
$ws.Range("B158").Value = 26
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 5
$ws.Range("E158").Value = 20
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 1
$ws.Range("A159").Value = "Haiti"
$ws.Range("B159").Value = 24
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 23
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 1
$ws.Range("A160").Value = "Tanzania"
$ws.Range("B160").Value = 24
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 5
$ws.Range("E160").Value = 18
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 1
$ws.Range("A161").Value = "Birmania"
$ws.Range("B161").Value = 22
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 21
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 1
$ws.Range("A162").Value = "Libia"
$ws.Range("B162").Value = 19
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 1
$ws.Range("E162").Value = 17
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 1
$ws.Range("A163").Value = "Siria"
$ws.Range("B163").Value = 19
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 2
$ws.Range("E163").Value = 15
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 2
$ws.Range("A164").Value = "Maldivas"
$ws.Range("B164").Value = 19
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 13
$ws.Range("E164").Value = 6
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0
$ws.Range("A165").Value = "Nueva Caledonia"
$ws.Range("B165").Value = 18
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 1
$ws.Range("E165").Value = 17
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0
$ws.Range("A166").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B166").Value = 17
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 17
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
$ws.Range("A167").Value = "Namibia"
$ws.Range("B167").Value = 16
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 3
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0
$ws.Range("A168").Value = "Guinea Ecuatorial"
$ws.Range("B168").Value = 16
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 3
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0
$ws.Range("A169").Value = "Angola"
$ws.Range("B169").Value = 16
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 2
$ws.Range("E169").Value = 12
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 2
$ws.Range("A170").Value = "Fiyi"
$ws.Range("B170").Value = 15
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 15
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0
$ws.Range("A171").Value = "Antigua y Barbuda"
$ws.Range("B171").Value = 15
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 0
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 1
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0
$ws.Range("A172").Value = "Dominica"
$ws.Range("B172").Value = 15
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 1
$ws.Range("E172").Value = 14
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0
$ws.Range("A173").Value = "Mongolia"
$ws.Range("B173").Value = 15
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 4
$ws.Range("E173").Value = 11
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
$ws.Range("A174").Value = "Laos"
$ws.Range("B174").Value = 14
$ws.Range("C174").Value = 2
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 14
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0
$ws.Range("A175").Value = "Santa Lucia"
$ws.Range("B175").Value = 14
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 1
$ws.Range("E175").Value = 13
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0
$ws.Range("A176").Value = "Sudan"
$ws.Range("B176").Value = 14
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 2
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 2
$ws.Range("A177").Value = "Liberia"
$ws.Range("B177").Value = 14
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 8
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 3
$ws.Range("A178").Value = "Curazao"
$ws.Range("B178").Value = 13
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 5
$ws.Range("E178").Value = 7
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 1
$ws.Range("A179").Value = "Granada"
$ws.Range("B179").Value = 12
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 12
$ws.Range("F179").Value = 2
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0
$ws.Range("A180").Value = "Seychelles"
$ws.Range("B180").Value = 11
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 11
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0
$ws.Range("A181").Value = "Groenlandia"
$ws.Range("B181").Value = 11
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 4
$ws.Range("E181").Value = 7
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0
$ws.Range("A182").Value = "San Cristobal y Nieves"
$ws.Range("B182").Value = 10
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0
$ws.Range("A183").Value = "Surinam"
$ws.Range("B183").Value = 10
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 1
$ws.Range("A184").Value = "Zimbabue"
$ws.Range("B184").Value = 10
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 9
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 1
$ws.Range("A185").Value = "Mozambique"
$ws.Range("B185").Value = 10
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 1
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0
$ws.Range("A186").Value = "Suazilandia"
$ws.Range("B186").Value = 10
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 4
$ws.Range("E186").Value = 6
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0
$ws.Range("A187").Value = "Republica del Chad"
$ws.Range("B187").Value = 9
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 9
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0
$ws.Range("A188").Value = "Nepal"
$ws.Range("B188").Value = 9
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0
$ws.Range("A189").Value = "Montserrat"
$ws.Range("B189").Value = 9
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 7
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 2
$ws.Range("A190").Value = "Republica de Africa Central"
$ws.Range("B190").Value = 8
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0
$ws.Range("A191").Value = "Malaui"
$ws.Range("B191").Value = 8
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 7
$ws.Range("F191").Value = 1
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 1
$ws.Range("A192").Value = "Islas Turcas y Caicos"
$ws.Range("B192").Value = 8
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 7
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 1
$ws.Range("A193").Value = "Santa Sede"
$ws.Range("B193").Value = 7
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 0
$ws.Range("E193").Value = 7
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0
$ws.Range("A194").Value = "San Vicente y las Granadinas"
$ws.Range("B194").Value = 7
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 1
$ws.Range("E194").Value = 6
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0
$ws.Range("A195").Value = "Somalia"
$ws.Range("B195").Value = 7
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 6
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0
$ws.Range("A196").Value = "Belice"
$ws.Range("B196").Value = 7
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 6
$ws.Range("F196").Value = 1
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1
$ws.Range("A197").Value = "Cabo Verde"
$ws.Range("B197").Value = 7
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 1
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 1
$ws.Range("A198").Value = "Sierra Leona"
$ws.Range("B198").Value = 6
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 6
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0
$ws.Range("A199").Value = "Nicaragua"
$ws.Range("B199").Value = 6
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1
$ws.Range("A200").Value = "Botsuana"
$ws.Range("B200").Value = 6
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 0
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1
$ws.Range("A201").Value = "San Bartolome"
$ws.Range("B201").Value = 6
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 1
$ws.Range("E201").Value = 5
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0
$ws.Range("A202").Value = "Mauritania"
$ws.Range("B202").Value = 6
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 3
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 1
$ws.Range("A203").Value = "Butan"
$ws.Range("B203").Value = 5
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 3
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0
$ws.Range("A204").Value = "Sahara Occidental"
$ws.Range("B204").Value = 4
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 0
$ws.Range("E204").Value = 4
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0
$ws.Range("A205").Value = "Santo Tome y Principe"
$ws.Range("B205").Value = 4
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("A206").Value = "Gambia"
$ws.Range("B206").Value = 4
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 2
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1
$ws.Range("A207").Value = "Burundi"
$ws.Range("B207").Value = 3
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 3
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0
$ws.Range("A208").Value = "Anguila"
$ws.Range("B208").Value = 3
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 3
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
$ws.Range("A209").Value = "Islas Virgenes Britanicas"
$ws.Range("B209").Value = 3
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 3
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 2
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 2
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
$ws.Range("A211").Value = "Papua Nueva Guinea"
$ws.Range("B211").Value = 2
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 0
$ws.Range("E211").Value = 2
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B212").Value = 2
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 2
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0
$ws.Range("A213").Value = "Timor Oriental"
$ws.Range("B213").Value = 1
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 0
$ws.Range("E213").Value = 1
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0
$ws.Range("A214").Value = "Sudan del Sur"
$ws.Range("B214").Value = 1
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 0
$ws.Range("E214").Value = 1
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
$ws.Range("A215").Value = "San Pedro y Miquelon"
$ws.Range("B215").Value = 1
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 1
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
